$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44308, 0, 3, 140.1214385801028),
    @(44309, 0, 3, 140.1214385801028),
    @(44310, 0, 2, 93.41429238673517),
    @(44311, 0, 0, 0),
    @(44312, 1, 1, 46.70714619336758)
)

$startRow = 234
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
